# Add a new worksheet "ODI Batting Extra" after the existing "ODI Bowling"
# sheet, populated with MATCH_CODE / BATTING_POSITION / NUM_4 / NUM_6 /
# PERCENT_RUNS_OF_TOTAL / MAN_OF_MATCH data, matching the author's commit.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ODI Batting Extra"

# ---- Header row (bold, centered, thin-bordered - matches the other tabs) ----
$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# ---- Data rows ----
# Row 2
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "4166"
$ws.Range("B2").Value = 8
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "0"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "0"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2.80%"
$ws.Range("F2").Value = "NO"

# Row 3
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "4169"
$ws.Range("F3").Value = "NO"

$ws.Range("A1").Select()
